$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the typo "inagenes" -> "imagenes" in the Sprint 1 feature text.
$ws.Range("B2").Value = "Captura de imagenes y procesamiento inicial de imagenes"

# Reword the Sprint 5 feature text.
$ws.Range("B15").Value = "relacion entre la interaccion de video con acciones de vision"

# Move the visible selection down to where the user was last working (B26:B27),
# scrolling the window so row 7 is at the top.
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("B26:B27").Select()
